{"js": "// Wrap the literal text \"A7C8CA\" with four bookmarks (OLE_LINK1, OLE_LINK2,\n// OLE_LINK3, OLE_LINK8) \u2014 these are the classic \"pasted from Office OLE\n// clipboard\" markers. This mirrors the diff: bookmarkStart elements are\n// inserted right before the run containing \"A7C8CA\" and the matching\n// bookmarkEnd elements right after it.\nconst body = context.document.body;\nconst results = body.search(\"A7C8CA\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Text \"A7C8CA\" not found in document body.');\n}\n\nconst range = results.items[0];\n// Office.js inserts bookmarks in nested fashion around the whole range, so\n// calling insertBookmark multiple times on the same range stacks the starts\n// immediately before the text and the ends immediately after it \u2014 exactly\n// the bookmarkStart x4 / bookmarkEnd x4 pairing in the target diff.\nrange.insertBookmark(\"OLE_LINK1\");\nrange.insertBookmark(\"OLE_LINK2\");\nrange.insertBookmark(\"OLE_LINK3\");\nrange.insertBookmark(\"OLE_LINK8\");\n\nawait context.sync();\n", "ps1": "# Wrap the literal text \"A7C8CA\" with four bookmarks (OLE_LINK1, OLE_LINK2,\n# OLE_LINK3, OLE_LINK8) \u2014 these are the classic \"pasted from Office OLE\n# clipboard\" markers. This mirrors the diff: bookmarkStart elements are\n# inserted right before the run containing \"A7C8CA\" and the matching\n# bookmarkEnd elements right after it.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$found = $rng.Find.Execute(\"A7C8CA\")\n\nif (-not $found) {\n    throw \"Text 'A7C8CA' not found in document content.\"\n}\n\n# Find.Execute collapses/extends $rng to the matched text in place, so the\n# same range is reused for all four Bookmarks.Add calls; Word nests the\n# bookmark starts immediately before the text and the ends immediately\n# after it, matching the target bookmarkStart x4 / bookmarkEnd x4 pairing.\n$d.Bookmarks.Add(\"OLE_LINK1\", $rng)\n$d.Bookmarks.Add(\"OLE_LINK2\", $rng)\n$d.Bookmarks.Add(\"OLE_LINK3\", $rng)\n$d.Bookmarks.Add(\"OLE_LINK8\", $rng)\n"}
